$d = $word.ActiveDocument

# 1. The lone "_GoBack" bookmark currently sits around the drawing in the
#    "Game Layout" section; it needs to move to the end of the intro
#    paragraph. Remove it from its old spot first.
$d.Bookmarks.Item("_GoBack").Delete()

# 2. Locate the intro paragraph's run precisely via Find (robust to any
#    offset drift) and capture its Start/End.
$rng = $d.Content
$found = $rng.Find.Execute(
    "A Platformer Game where the Player controls a Rocket in space by tilting the device side to side, avoiding obstacles and trying to keep the Rocket moving.",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the intro paragraph text to update."
}

# Re-anchor into a fresh Range object using the Find hit's coordinates;
# InsertXML misbehaves when called directly on a Range still owned by Find.
$target = $d.Range($rng.Start, $rng.End)

# 3. Replace that run with the new two-run split plus a collapsed
#    "_GoBack" bookmark placed right after the new text, using InsertXML
#    so the exact run/bookmark structure is produced.
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body><w:p>' +
       '<w:r><w:t xml:space="preserve">A Platformer Game where the Player controls a Rocket in space by tilting the device side </w:t></w:r>' +
       '<w:r><w:t>to side, avoiding obstacles, collecting Star Shards and accumulating Score.</w:t></w:r>' +
       '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
       '</w:p></w:body></w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xml)
